$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values (columns M through T) with new TPM-derived numbers.
$ws.Range("M2").Value = 0.044174
$ws.Range("N2").Value = 0.132522
$ws.Range("O2").Value = 0.2474193313505733
$ws.Range("P2").Value = 0.2474193313505733
$ws.Range("Q2").Value = 0.006964988203333334
$ws.Range("R2").Value = 0.06268489383
$ws.Range("S2").Value = 0.2474193313505733
$ws.Range("T2").Value = 0.2474193313505733

# Row 2's Target cluster label changes from "MuSCs" to "ECs" as a new
# cluster label is introduced ahead of it in the shared-string table.
$ws.Range("D2").Value = "ECs"

# Add a new data row (row 3) for the Rspo1 -> Lgr6 interaction targeting MuSCs.
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Rspo1"
$ws.Range("C3").Value = "Lgr6"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1576716666666667
$ws.Range("H3").Value = 0.473015
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.134365
$ws.Range("N3").Value = 0.403095
$ws.Range("O3").Value = 0.7525806686494267
$ws.Range("P3").Value = 0.7525806686494266
$ws.Range("Q3").Value = 0.02118555349166667
$ws.Range("R3").Value = 0.190669981425
$ws.Range("S3").Value = 0.7525806686494267
$ws.Range("T3").Value = 0.7525806686494266
